# Generate Report for Handback
# For each language sheet (zh-cn, de-de), the handoff rows (row 2 and row 3)
# are updated to reflect a completed handback:
#   - Status (col B) becomes "Handed back: in sync with en-US"
#   - Latest Target File (col E) / Latest Handback File (col F) are filled in
#     (same filenames as the Source File Name / Latest Handoff File columns)
#   - Latest Handback DateTime (col G) gets a real timestamp instead of the
#     "0001-01-01 00:00:00" placeholder

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$sheetInfo = @(
    @{
        Name = "zh-cn"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2036d52ade0a2f96bd2f92f92f75e722b966832f/e2e"
        XlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/aa54a4a6c8477c755fb6bdf44d7f05fd37b89ad6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
        ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2036d52ade0a2f96bd2f92f92f75e722b966832f/.localization-config"
        Row2Xlf = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.5ae0b5988e282534b84e59c077b189920b473614.zh-cn.xlf"
        Row3Xlf = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.772bc25e666e30d138fa3394e3ca21f22c886e94.zh-cn.xlf"
        Row2HandbackTime = "2016-03-08 19:11:05"
        Row3HandbackTime = "2016-03-08 19:11:05"
    },
    @{
        Name = "de-de"
        MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2036d52ade0a2f96bd2f92f92f75e722b966832f/e2e"
        XlfUrlBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c17c68372f32aa4ef7b32b75b9f8c2667ec1d6db/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
        ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/2036d52ade0a2f96bd2f92f92f75e722b966832f/.localization-config"
        Row2Xlf = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.5ae0b5988e282534b84e59c077b189920b473614.de-de.xlf"
        Row3Xlf = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.772bc25e666e30d138fa3394e3ca21f22c886e94.de-de.xlf"
        Row2HandbackTime = "2016-03-08 19:11:30"
        Row3HandbackTime = "2016-03-08 19:11:30"
    }
)

$mdRow2 = "7712ca7f-55cb-4ea6-8b16-8a3a3487afa2.md"
$mdRow3 = "b07d73ed-43ea-4d82-90ed-f4f09c252c2b.md"
$configName = ".localization-config"

# The "Ready for handoff" status string is shared (by shared-string index)
# across the Overview sheet AND each language sheet, so replace it everywhere
# in one pass first.
foreach ($ws0 in $wb.Worksheets) {
    $ws0.Cells.Replace("Ready for handoff", $statusHandedBack)
}

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # --- Fill in Latest Target File / Latest Handback File text values ---
    $ws.Range("E2").Value = $mdRow2
    $ws.Range("F2").Value = $info.Row2Xlf
    $ws.Range("E3").Value = $mdRow3
    $ws.Range("F3").Value = $info.Row3Xlf

    # --- Latest Handback DateTime now has a real timestamp ---
    $ws.Range("G2").Value = $info.Row2HandbackTime
    $ws.Range("G3").Value = $info.Row3HandbackTime

    # --- Rebuild hyperlinks in row order (A2,C2,E2,F2,A3,C3,E3,F3,A4) so the
    #     relationship ids line up the same way Excel would renumber them. ---
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), "$($info.MdUrl)/$mdRow2", "", "", $mdRow2)
    $ws.Hyperlinks.Add($ws.Range("C2"), "$($info.XlfUrlBase)/$($info.Row2Xlf)", "", "", $info.Row2Xlf)
    $ws.Hyperlinks.Add($ws.Range("E2"), "$($info.MdUrl)/$mdRow2", "", "", $mdRow2)
    $ws.Hyperlinks.Add($ws.Range("F2"), "$($info.XlfUrlBase)/$($info.Row2Xlf)", "", "", $info.Row2Xlf)

    $ws.Hyperlinks.Add($ws.Range("A3"), "$($info.MdUrl)/$mdRow3", "", "", $mdRow3)
    $ws.Hyperlinks.Add($ws.Range("C3"), "$($info.XlfUrlBase)/$($info.Row3Xlf)", "", "", $info.Row3Xlf)
    $ws.Hyperlinks.Add($ws.Range("E3"), "$($info.MdUrl)/$mdRow3", "", "", $mdRow3)
    $ws.Hyperlinks.Add($ws.Range("F3"), "$($info.XlfUrlBase)/$($info.Row3Xlf)", "", "", $info.Row3Xlf)

    $ws.Hyperlinks.Add($ws.Range("A4"), $info.ConfigUrl, "", "", $configName)

    # --- Match the existing hyperlink look (underline + blue) used by the
    #     other link cells in this table instead of Excel's default theme link. ---
    foreach ($addr in @("A2", "C2", "E2", "F2", "A3", "C3", "E3", "F3", "A4")) {
        $ws.Range($addr).Font.Underline = 2
        $ws.Range($addr).Font.Color = 15570276
    }
}
